$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1, matching style of existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for columns I (I0) and J (IF), rows 2-81
$colI = @(5, 7, 8, 8, 8, 8, 7, 7, 8, 8, 8, 7, 6, 8, 7, 8, 7, 8, 8, 8, 8, 7, 8, 8, 8, 8, 8, 7, 7, 6, 9, 9, 8, 8, 9, 8, 9, 9, 9, 9, 9, 7, 9, 8, 8, 9, 8, 8, 8, 8, 8, 8, 10, 9, 6, 8, 8, 8, 8, 8, 7, 8, 7, 7, 8, 8, 8, 8, 8, 7, 6, 6, 6, 6, 7, 7, 6, 8, 6, 7)
$colJ = @(6, 7, 8, 8, 8, 8, 7, 8, 8, 8, 8, 8, 7, 8, 7, 8, 7, 8, 8, 8, 8, 7, 9, 8, 8, 8, 8, 7, 7, 7, 9, 9, 8, 8, 9, 8, 9, 9, 9, 9, 9, 7, 9, 8, 8, 9, 8, 8, 8, 8, 8, 8, 10, 9, 6, 8, 8, 8, 8, 8, 8, 8, 7, 7, 8, 8, 8, 8, 8, 7, 7, 6, 6, 6, 7, 7, 6, 8, 6, 7)

for ($i = 0; $i -lt $colI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $colI[$i]
    $ws.Cells.Item($row, 10).Value = $colJ[$i]
}
